$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colLetters = @("B","C","D","E","F","H","I","M")
$colIndex = @{ "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "H"=8; "I"=9; "M"=13 }

$data = @{
    2 = @(0.8592388859344169, 0.3418508317705005, 0.03885608809694219, 0.08718440659772142, 1.557330138999845, 0.07973214163530429, 0.9864072933334995, 0.3925723809337072)
    3 = @(0.7684864284123023, 0.3006474811034821, 0.03843379391268087, 0.08181199126970284, 1.488194764439896, 0.07973214163530429, 0.9564341054137628, 0.3545240259354401)
    4 = @(0.7131505995820078, 0.2754571263163825, 0.03819084857078536, 0.07857295233067418, 1.4467161335497, 0.07973214163530429, 0.938591943007097, 0.3313776483317739)
    5 = @(0.6906967838333742, 0.2652181449354316, 0.03809587335076259, 0.07726771936905763, 1.430054052876955, 0.07973214163530429, 0.9314604687329364, 0.321998475583392)
    6 = @(0.6869741073578268, 0.2635195280580831, 0.03808034358172563, 0.07705186575586254, 1.427301774360586, 0.07973214163530429, 0.9302846556909117, 0.3204442527094216)
    7 = @(0.7128473925133676, 0.2753189347175748, 0.03818955150664394, 0.07855529036411824, 1.446490452300353, 0.07973214163530429, 0.9384952036657239, 0.3312509435946254)
    8 = @(0.8278664377279483, 0.3276207234193294, 0.03870703842743239, 0.08531943730441327, 1.533288758455726, 0.07973214163530429, 0.9759548793168165, 0.3794079595848103)
    9 = @(1.056554149964256, 0.4310988719786906, 0.03985532292296057, 0.09907132646731043, 1.711361868640807, 0.07973214163530429, 1.053955803326403, 0.4756047487398121)
    10 = @(1.226597565596251, 0.5077623236032878, 0.04078569740407545, 0.109492700764342, 1.847220030664943, 0.07973214163530429, 1.114159691422955, 0.5474332886716837)
    11 = @(1.304421482805367, 0.5427952289990685, 0.04122894682303979, 0.1143072629545614, 1.910169149029656, 0.07973214163530429, 1.142204973610276, 0.5803781312048102)
    12 = @(1.333960866045061, 0.5560853916527435, 0.04139976655374511, 0.1161413861724441, 1.934175194312019, 0.07973214163530429, 1.152921812438493, 0.5928935943052096)
    13 = @(1.327595930807604, 0.5532220323208321, 0.04136284386292743, 0.1157458823683797, 1.928997505798833, 0.07973214163530429, 1.15060941593299, 0.5901963687232836)
    14 = @(1.306850310752736, 0.5438881315085951, 0.04124294019405994, 0.1144579359140394, 1.912140740838993, 0.07973214163530429, 1.143084704606522, 0.5814069780152806)
    15 = @(1.294152071615372, 0.5381740042017782, 0.04116988537495558, 0.1136704680071929, 1.901837550442167, 0.07973214163530429, 1.138488255788502, 0.5760284628549641)
    16 = @(1.221521182559854, 0.5054761282407867, 0.04075714068621039, 0.1091795720662248, 1.843129545314468, 0.07973214163530429, 1.11234028605476, 0.5452857932367863)
    17 = @(1.177086010770552, 0.4854584321332709, 0.04050912755196379, 0.1064436997860128, 1.807410370923662, 0.07973214163530429, 1.096469392312216, 0.5264960689365097)
    18 = @(1.151572248965067, 0.4739596112447657, 0.04036835552144424, 0.1048770269032104, 1.786973341977529, 0.07973214163530429, 1.08740270639295, 0.5157140274994703)
    19 = @(1.142941268937648, 0.4700688188489721, 0.04032101237529417, 0.1043477581770205, 1.78007209966853, 0.07973214163530429, 1.084343440057268, 0.5120677249783796)
    20 = @(1.181811627981801, 0.4875878067832105, 0.04053533388334785, 0.1067342185639149, 1.811201565235393, 0.07973214163530429, 1.098152461867727, 0.5284936390555686)
    21 = @(1.31294191298548, 0.5466290644652076, 0.04127807744130507, 0.1148359370380732, 1.917087376346899, 0.07973214163530429, 1.14529225393305, 0.5839875408696003)
    22 = @(1.399046968058883, 0.585356179735868, 0.04178086281141447, 0.1201948807503115, 1.987274241148214, 0.07973214163530429, 1.176665433310617, 0.6204895364001572)
    23 = @(1.353053611058556, 0.5646735557010629, 0.04151089767735527, 0.1173287428476613, 1.949722884996618, 0.07973214163530429, 1.159868634851733, 0.6009859683627354)
    24 = @(1.179675074980764, 0.4866250866321593, 0.0405234803574146, 0.1066028556867167, 1.809487260072018, 0.07973214163530429, 1.097391366546461, 0.5275904740860824)
    25 = @(0.9943406728968966, 0.4029991464116733, 0.03952983979668545, 0.09529676884171323, 1.662322545892408, 0.07973214163530429, 1.032355068990114, 0.4493843379078157)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($j = 0; $j -lt $colLetters.Length; $j++) {
        $colLetter = $colLetters[$j]
        $colNum = $colIndex[$colLetter]
        $ws.Cells.Item($row, $colNum).Value = $vals[$j]
    }
}
